$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename sheet / tab: "CYRS review" -> "SRS review"
# ---------------------------------------------------------------------------
$ws.Name = "SRS review"

# ---------------------------------------------------------------------------
# 2. Close review points 2-6 (column F: Open -> Closed)
# ---------------------------------------------------------------------------
$ws.Range("F2:F6").Value = "Closed"

# Row heights grow a bit once re-saved/re-wrapped in the newer file
$ws.Range("A2").RowHeight = 30
$ws.Range("A3").RowHeight = 75
$ws.Range("A4").RowHeight = 60
$ws.Range("A5").RowHeight = 45
$ws.Range("A6").RowHeight = 75

# ---------------------------------------------------------------------------
# 3. Add two new (open) review points in rows 7 and 8, matching the layout
#    and formatting already used by row 6 (copy formats, then overwrite the
#    content of the new rows).
# ---------------------------------------------------------------------------
$ws.Range("A6:F6").Copy() | Out-Null
$ws.Range("A7:F8").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A7").Value = "13/2/2020"
$ws.Range("B7").Value = "Ali"
$ws.Range("C7").Value = "SRS"
$ws.Range("D7").Value = "Tones that will be sent to buzzer aren't defined, For example the tone signal is DC signal or PWM and if it's PWM what's the duty and freq. for it?"
$ws.Range("E7").Value = ""
$ws.Range("F7").Value = "Open"
$ws.Range("A7").RowHeight = 45

$ws.Range("A8").Value = "13/2/2020"
$ws.Range("B8").Value = "Ali"
$ws.Range("C8").Value = "SRS"
$ws.Range("D8").Value = "Cur_x and Cur_y are not deined with specific values, developer shall know the exact values for them"
$ws.Range("E8").Value = ""
$ws.Range("F8").Value = "Open"
$ws.Range("A8").RowHeight = 30

# ---------------------------------------------------------------------------
# 4. Conditional formatting for the new rows, mirroring the rules already
#    applied to F2:F6 / E2:E6 (Open/Closed text highlighting, Accepted/
#    Rejected text highlighting).
# ---------------------------------------------------------------------------
# xlTextString = 9, xlContains = 0
$f7a = $ws.Range("F7").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Open")
$f7a.Font.Color = 255
$f7b = $ws.Range("F7").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Closed")
$f7b.Font.Color = 24832
$f7b.Interior.Color = 13561798
$f7c = $ws.Range("F7").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Open")
$f7c.Font.Color = 393372

# xlCellValue = 1, xlEqual = 3
$e7a = $ws.Range("E7").FormatConditions.Add(1, 3, '"Rejected"')
$e7a.Font.Color = 393372
$e7b = $ws.Range("E7").FormatConditions.Add(1, 3, '"Accepted"')
$e7b.Font.Color = 24832
$e7b.Interior.Color = 13561798

$f8a = $ws.Range("F8").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Open")
$f8a.Font.Color = 255
$f8b = $ws.Range("F8").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Closed")
$f8b.Font.Color = 24832
$f8b.Interior.Color = 13561798
$f8c = $ws.Range("F8").FormatConditions.Add(9, 0, [System.Type]::Missing, [System.Type]::Missing, "Open")
$f8c.Font.Color = 393372

$e8a = $ws.Range("E8").FormatConditions.Add(1, 3, '"Rejected"')
$e8a.Font.Color = 393372
$e8b = $ws.Range("E8").FormatConditions.Add(1, 3, '"Accepted"')
$e8b.Font.Color = 24832
$e8b.Interior.Color = 13561798

# ---------------------------------------------------------------------------
# 5. Extend the data-validation drop-downs to cover the new rows too.
# ---------------------------------------------------------------------------
$ws.Range("F2:F6").Validation.Delete()
$ws.Range("E2:E6").Validation.Delete()
$ws.Range("F2:F8").Validation.Add(3, 1, 1, '"Open, Closed"')
$ws.Range("E2:E8").Validation.Add(3, 1, 1, '"Accepted, Rejected"')

# ---------------------------------------------------------------------------
# 6. Update the view: scroll down a little and select D7 (the first of the
#    newly-added review points) as the active cell.
# ---------------------------------------------------------------------------
$ws.Range("D7").Select() | Out-Null
